$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of the per-observation data across
# rows 2, 3 and 5 (header is row 1, row 4 and row 6 are untouched):
#   new row 2 <- old row 3
#   new row 3 <- old row 5
#   new row 5 <- old row 2
# Columns A, B, D, E, F, G, P, Q, R, AI carry the row-specific data that
# moves; C, H, I and the rest of the row stay identical across the three
# rows so they don't need to be touched.

$cols = @("A", "B", "D", "E", "F", "G", "P", "Q", "R", "AI")
$rows = @(2, 3, 5)

# Snapshot the current ("before") values for the rows involved.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping of destination row -> source row (cyclic rotation).
$sourceFor = @{ 2 = 3; 3 = 5; 5 = 2 }

foreach ($destRow in $rows) {
    $srcRow = $sourceFor[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $val = $srcData[$col]
        if ($col -eq "AI") {
            if ($val -eq $null) {
                $ws.Range("AI$destRow").ClearContents()
            } else {
                $ws.Range("AI$destRow").Value = $val
            }
        } else {
            $ws.Range("$col$destRow").Value = $val
        }
    }
}
